# Update "想去人数" (wanted-to-go count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value  = 487
$ws1.Range("F6").Value  = 878
$ws1.Range("F9").Value  = 2148
$ws1.Range("F10").Value = 620
$ws1.Range("F11").Value = 281
$ws1.Range("F12").Value = 116
$ws1.Range("F13").Value = 1051
$ws1.Range("F14").Value = 176
$ws1.Range("F15").Value = 2177
$ws1.Range("F16").Value = 642
$ws1.Range("F17").Value = 11951
$ws1.Range("F18").Value = 1216
$ws1.Range("F19").Value = 552
$ws1.Range("F21").Value = 13
$ws1.Range("F23").Value = 36
$ws1.Range("F24").Value = 258
$ws1.Range("F27").Value = 15

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 79
$ws2.Range("F12").Value = 56

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5685

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 5685
$ws4.Range("F9").Value  = 487
$ws4.Range("F10").Value = 878
$ws4.Range("F14").Value = 2148
$ws4.Range("F15").Value = 620
$ws4.Range("F16").Value = 281
$ws4.Range("F18").Value = 116
$ws4.Range("F20").Value = 1051
$ws4.Range("F22").Value = 176
$ws4.Range("F25").Value = 2177
$ws4.Range("F26").Value = 642
$ws4.Range("F27").Value = 79
$ws4.Range("F28").Value = 56
$ws4.Range("F29").Value = 1216
$ws4.Range("F30").Value = 552
$ws4.Range("F32").Value = 13
$ws4.Range("F35").Value = 36
$ws4.Range("F38").Value = 258
$ws4.Range("F49").Value = 15
